# Update the 20x5 arithmetic-practice table in place.
# Each cell holds one "a OP b = c" expression; replace by (row, col)
# position rather than by text-match, since several expressions repeat
# (e.g. "92-78=14" occurs twice but maps to two different replacements).
$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text = "50-1=49"
$tbl.Cell(1,2).Range.Text = "52-30=22"
$tbl.Cell(1,3).Range.Text = "0+71=71"
$tbl.Cell(1,4).Range.Text = "32+44=76"
$tbl.Cell(1,5).Range.Text = "94-79=15"
$tbl.Cell(2,1).Range.Text = "15+40=55"
$tbl.Cell(2,2).Range.Text = "26+59=85"
$tbl.Cell(2,3).Range.Text = "9-8=1"
$tbl.Cell(2,4).Range.Text = "99-49=50"
$tbl.Cell(2,5).Range.Text = "49+40=89"
$tbl.Cell(3,1).Range.Text = "17+41=58"
$tbl.Cell(3,2).Range.Text = "11+69=80"
$tbl.Cell(3,3).Range.Text = "70-36=34"
$tbl.Cell(3,4).Range.Text = "68-56=12"
$tbl.Cell(3,5).Range.Text = "55-52=3"
$tbl.Cell(4,1).Range.Text = "68-29=39"
$tbl.Cell(4,2).Range.Text = "75-25=50"
$tbl.Cell(4,3).Range.Text = "65+9=74"
$tbl.Cell(4,4).Range.Text = "51-14=37"
$tbl.Cell(4,5).Range.Text = "16-3=13"
$tbl.Cell(5,1).Range.Text = "96-14=82"
$tbl.Cell(5,2).Range.Text = "56+35=91"
$tbl.Cell(5,3).Range.Text = "67-15=52"
$tbl.Cell(5,4).Range.Text = "51-9=42"
$tbl.Cell(5,5).Range.Text = "7+26=33"
$tbl.Cell(6,1).Range.Text = "63-50=13"
$tbl.Cell(6,2).Range.Text = "61-47=14"
$tbl.Cell(6,3).Range.Text = "88-65=23"
$tbl.Cell(6,4).Range.Text = "38+42=80"
$tbl.Cell(6,5).Range.Text = "44+27=71"
$tbl.Cell(7,1).Range.Text = "87-76=11"
$tbl.Cell(7,2).Range.Text = "40+57=97"
$tbl.Cell(7,3).Range.Text = "42-7=35"
$tbl.Cell(7,4).Range.Text = "8+82=90"
$tbl.Cell(7,5).Range.Text = "43+34=77"
$tbl.Cell(8,1).Range.Text = "85+9=94"
$tbl.Cell(8,2).Range.Text = "43+55=98"
$tbl.Cell(8,3).Range.Text = "63-5=58"
$tbl.Cell(8,4).Range.Text = "57+28=85"
$tbl.Cell(8,5).Range.Text = "87-80=7"
$tbl.Cell(9,1).Range.Text = "87-49=38"
$tbl.Cell(9,2).Range.Text = "46-15=31"
$tbl.Cell(9,3).Range.Text = "52+16=68"
$tbl.Cell(9,4).Range.Text = "99-74=25"
$tbl.Cell(9,5).Range.Text = "40+5=45"
$tbl.Cell(10,1).Range.Text = "62-14=48"
$tbl.Cell(10,2).Range.Text = "66-12=54"
$tbl.Cell(10,3).Range.Text = "98-97=1"
$tbl.Cell(10,4).Range.Text = "67-64=3"
$tbl.Cell(10,5).Range.Text = "68-63=5"
$tbl.Cell(11,1).Range.Text = "91-2=89"
$tbl.Cell(11,2).Range.Text = "35-21=14"
$tbl.Cell(11,3).Range.Text = "5+54=59"
$tbl.Cell(11,4).Range.Text = "15-9=6"
$tbl.Cell(11,5).Range.Text = "3+27=30"
$tbl.Cell(12,1).Range.Text = "72+6=78"
$tbl.Cell(12,2).Range.Text = "12+69=81"
$tbl.Cell(12,3).Range.Text = "72-22=50"
$tbl.Cell(12,4).Range.Text = "90-52=38"
$tbl.Cell(12,5).Range.Text = "30-12=18"
$tbl.Cell(13,1).Range.Text = "61-56=5"
$tbl.Cell(13,2).Range.Text = "87-6=81"
$tbl.Cell(13,3).Range.Text = "41-21=20"
$tbl.Cell(13,4).Range.Text = "19-13=6"
$tbl.Cell(13,5).Range.Text = "70-52=18"
$tbl.Cell(14,1).Range.Text = "69-32=37"
$tbl.Cell(14,2).Range.Text = "58+38=96"
$tbl.Cell(14,3).Range.Text = "80+0=80"
$tbl.Cell(14,4).Range.Text = "15+64=79"
$tbl.Cell(14,5).Range.Text = "99-93=6"
$tbl.Cell(15,1).Range.Text = "30+10=40"
$tbl.Cell(15,2).Range.Text = "13+39=52"
$tbl.Cell(15,3).Range.Text = "28-21=7"
$tbl.Cell(15,4).Range.Text = "13+34=47"
$tbl.Cell(15,5).Range.Text = "85-28=57"
$tbl.Cell(16,1).Range.Text = "75-56=19"
$tbl.Cell(16,2).Range.Text = "74-13=61"
$tbl.Cell(16,3).Range.Text = "20-5=15"
$tbl.Cell(16,4).Range.Text = "95-41=54"
$tbl.Cell(16,5).Range.Text = "54-44=10"
$tbl.Cell(17,1).Range.Text = "44+45=89"
$tbl.Cell(17,2).Range.Text = "37-13=24"
$tbl.Cell(17,3).Range.Text = "34-18=16"
$tbl.Cell(17,4).Range.Text = "38+21=59"
$tbl.Cell(17,5).Range.Text = "92-88=4"
$tbl.Cell(18,1).Range.Text = "19+24=43"
$tbl.Cell(18,2).Range.Text = "49+17=66"
$tbl.Cell(18,3).Range.Text = "78-64=14"
$tbl.Cell(18,4).Range.Text = "69+29=98"
$tbl.Cell(18,5).Range.Text = "21-21=0"
$tbl.Cell(19,1).Range.Text = "58-5=53"
$tbl.Cell(19,2).Range.Text = "34+17=51"
$tbl.Cell(19,3).Range.Text = "99-13=86"
$tbl.Cell(19,4).Range.Text = "73+23=96"
$tbl.Cell(19,5).Range.Text = "31+44=75"
$tbl.Cell(20,1).Range.Text = "34+13=47"
$tbl.Cell(20,2).Range.Text = "48-43=5"
$tbl.Cell(20,3).Range.Text = "49+21=70"
$tbl.Cell(20,4).Range.Text = "58+6=64"
$tbl.Cell(20,5).Range.Text = "79-43=36"
